# Update the workbook to add carjacking data through 2021-09-24
# (previously through 2021-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "Through 2021-09-24"

# 2. Update the column header / shared string text for the "current" month
#    column (column B), which describes the cutoff date.
$ws.Range("B1").Value = "September 2021 (through September 24)"

# 3. Update the per-neighborhood counts for column B (current-month totals)
#    and other affected cells for 2021-10-02 data.

# Row 2 - Garfield Park
$ws.Range("AL2").Value = 6

# Row 3 - North Lawndale
$ws.Range("B3").Value = 9
$ws.Range("K3").Value = 8

# Row 4 - Humboldt Park
$ws.Range("B4").Value = 5

# Row 5 - Austin
$ws.Range("AU5").Value = 2

# Row 6 - Roseland
$ws.Range("AC6").Value = 4

# Row 11 - Kenwood
$ws.Range("T11").Value = 2

# Row 13 - Chatham
$ws.Range("D13").Value = 5

# Row 18 - Grand Boulevard
$ws.Range("B18").Value = 2

# Row 19 - Wicker Park
$ws.Range("B19").Value = 5
$ws.Range("K19").Value = 1

# Row 21 - River North
$ws.Range("AL21").Value = 2

# Row 23 - Ashburn
$ws.Range("K23").Value = 5

# Row 32 - Chicago Lawn
$ws.Range("B32").Value = 3

# Row 33 - Lake View
$ws.Range("B33").Value = 4

# Row 51 - Loop
$ws.Range("T51").Value = 2

# Row 55 - Grand Crossing
$ws.Range("K55").Value = 9

# Row 59 - Archer Heights
$ws.Range("B59").Value = 1

# Row 68 - East Side
$ws.Range("K68").Value = 1

# Row 71 - Galewood
$ws.Range("AL71").Value = 1

# Row 89 - Oakland
$ws.Range("T89").Value = 1
